$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 30. This shifts the existing rows 30-48
# (and their data) down to rows 31-49, preserving all their values.
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new weekly price record.
$ws.Cells.Item(30, 1).Value = 10
$ws.Cells.Item(30, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(30, 3).Value = "La Araucanía"
$ws.Cells.Item(30, 4).Value = 44839
$ws.Cells.Item(30, 5).Value = 9
$ws.Cells.Item(30, 6).Value = 100112042
$ws.Cells.Item(30, 7).Value = "Locoto"
$ws.Cells.Item(30, 8).Value = "Sin especificar"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 80
$ws.Cells.Item(30, 11).Value = 2500
$ws.Cells.Item(30, 12).Value = 2500
$ws.Cells.Item(30, 13).Value = 2500
$ws.Cells.Item(30, 14).Value = "`$/kilo"
$ws.Cells.Item(30, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(30, 16).Value = 2500
$ws.Cells.Item(30, 17).Value = 1
$ws.Cells.Item(30, 18).Value = "Hortaliza"
